# Change User Password Functionality
# Adds a "Change Password" section of locators to the "Web" sheet, right
# after the existing "Edit User" block, following the same table layout
# (ElementID | ElementPath | Method) with blank separator rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: existing rows 27-28 (EditUser_Action / EditUser_UpdateButton)
#    need to end up at rows 28-29, followed by a new EditUser_ToastMessage
#    row, a blank separator, the new ChangePassword block (5 rows), and a
#    trailing blank block (4 rows) to mirror the rest of the sheet.
# ---------------------------------------------------------------------

# Blank row inserted above the existing EditUser_Action/UpdateButton pair
$ws.Rows.Item(27).Insert()

# New row for the "toast message" locator right after EditUser_UpdateButton
$ws.Rows.Item(30).Insert()

# Blank separator row before the new Change Password section
$ws.Rows.Item(31).Insert()

# Four trailing blank rows at the very end of the sheet (37-40)
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).Insert()

# ---------------------------------------------------------------------
# 2) Fill in the new / shifted content.
#    A lone "'" forces an explicit (empty) text value on a cell instead of
#    leaving it absent, matching the blank-but-present cells used
#    throughout this locator table as row separators.
# ---------------------------------------------------------------------

$ws.Range("A27").Value() = "'"
$ws.Range("B27").Value() = "'"
$ws.Range("C27").Value() = "'"

$ws.Range("A30").Value() = "EditUser_ToastMessage"
$ws.Range("B30").Value() = "'"
$ws.Range("C30").Value() = "By.xpath"

$ws.Range("A31").Value() = "'"
$ws.Range("B31").Value() = "'"
$ws.Range("C31").Value() = "'"

$ws.Range("A32").Value() = "UserChangePassword_Action"
$ws.Range("B32").Value() = '//*[@id="root"]/div[1]/div/div[2]/div/div[2]/div/div[3]/div/table/tbody/tr[3]/td[6]/div/button[2]'
$ws.Range("C32").Value() = "By.xpath"

$ws.Range("A33").Value() = "User_NewPassword"
$ws.Range("B33").Value() = "(//input[@placeholder='xxxxxxxx'])[1]"
$ws.Range("C33").Value() = "By.xpath"

$ws.Range("A34").Value() = "User_ConfirmPassword"
$ws.Range("B34").Value() = "(//input[@placeholder='xxxxxxxx'])[2]"
$ws.Range("C34").Value() = "By.xpath"

$ws.Range("A35").Value() = "User_PasswordResetButton"
$ws.Range("B35").Value() = "//button[normalize-space()='Reset']"
$ws.Range("C35").Value() = "By.xpath"

$ws.Range("A36").Value() = "ChangePassword_ToastMessage"
$ws.Range("B36").Value() = "'"
$ws.Range("C36").Value() = "By.xpath"

$ws.Range("A37").Value() = "'"
$ws.Range("B37").Value() = "'"
$ws.Range("C37").Value() = "'"

$ws.Range("A38").Value() = "'"
$ws.Range("B38").Value() = "'"
$ws.Range("C38").Value() = "'"

$ws.Range("A39").Value() = "'"
$ws.Range("B39").Value() = "'"
$ws.Range("C39").Value() = "'"

$ws.Range("A40").Value() = "'"
$ws.Range("B40").Value() = "'"
$ws.Range("C40").Value() = "'"

Write-Host "Change User Password Functionality: locators added (rows 27-40)."
